$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# EMU -> points conversion (1 pt = 12700 EMU), AddTextbox expects points.
$left   = 6344066 / 12700
$top    = 5775909 / 12700
$width  = 5225097 / 12700
$height = 400110 / 12700

$tb = $s.Shapes.AddTextbox(1, $left, $top, $width, $height)
$tb.Name = "CuadroTexto 5"

$tb.Fill.Visible = $false

$tf = $tb.TextFrame
$tf.Orientation = 1
$tf.WordWrap = $true
$tf.AutoSize = 1
$tf.VerticalAnchor = 3
$tf.HorizontalAnchor = $false

$tr = $tf.TextRange
$tr.Text = "JSON:"
$tr.Font.Size = 10
$tr.LanguageID = "es-MX"

$r2 = $tr.InsertAfter("`rmi-carrito")
$r2.Font.Size = 10
$r2.LanguageID = "es-MX"
